# Character-info sheet: keep only the header row and the single data row
# for "段总" / "生锈的斩牛刀" / "伊森利恩" / "盗贼" (originally row 7), discarding
# every other character row. Row/outline metadata (sheetFormatPr
# outlineLevelRow/outlineLevelCol) and the selected cell are preserved to
# match the target workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the outline-level bookkeeping on rows (1) and columns (3) that
# the original sheet carried in sheetFormatPr, independent of which rows
# end up staying in sheetData.
$ws.Outline.ShowLevels(1, 3)

# Delete the unwanted rows from the bottom up so earlier row numbers keep
# referring to the same original rows while we work.
$ws.Rows("8:18").Delete()
$ws.Rows("2:6").Delete()

# Restore the selection to match the target view state.
$ws.Range("D12").Select()

Write-Host "Trimmed character table to header + 段总 row"
